$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns are treated as text so values such as
# "37.00", "0.990" or "1.00" are not silently coerced into numbers
# and lose their original formatting/trailing zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '59.231.93'
$ws.Range("E2").Value = '  +0.70%  '

$ws.Range("D3").Value = '2.689.43'
$ws.Range("E3").Value = '  +5.25%  '

$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").Value = '519.06'
$ws.Range("E5").Value = '  +2.83%  '

$ws.Range("D6").Value = '146.21'
$ws.Range("E6").Value = '  +3.16%  '

$ws.Range("D7").Value = '0.994'
$ws.Range("E7").Value = '  -0.54%  '

$ws.Range("D8").Value = '0.570'
$ws.Range("E8").Value = '  +3.35%  '

$ws.Range("D9").Value = '2.731.29'
$ws.Range("E9").Value = '  +6.84%  '

$ws.Range("D10").Value = '6.26'
$ws.Range("E10").Value = '  +0.98%  '

$ws.Range("D11").Value = '0.108'
$ws.Range("E11").Value = '  +6.97%  '

$ws.Range("D12").Value = '0.338'
$ws.Range("E12").Value = '  +2.99%  '

$ws.Range("E13").Value = '  -0.93%  '

$ws.Range("D14").Value = '3.169.62'
$ws.Range("E14").Value = '  +5.63%  '

$ws.Range("D15").Value = '59.184.62'
$ws.Range("E15").Value = '  +0.66%  '

$ws.Range("D16").Value = '21.21'
$ws.Range("E16").Value = '  +3.62%  '

$ws.Range("D17").Value = '0.0000139'
$ws.Range("E17").Value = '  +3.71%  '

$ws.Range("D18").Value = '2.717.31'
$ws.Range("E18").Value = '  +6.29%  '

$ws.Range("D19").Value = '4.58'
$ws.Range("E19").Value = '  +1.99%  '

$ws.Range("D20").Value = '347.25'
$ws.Range("E20").Value = '  +5.15%  '

$ws.Range("D21").Value = '10.54'
$ws.Range("E21").Value = '  +5.35%  '

$ws.Range("D22").Value = '6.24'
$ws.Range("E22").Value = '  +5.67%  '

$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.25%  '

$ws.Range("D24").Value = '61.25'
$ws.Range("E24").Value = '  +3.19%  '

$ws.Range("D25").Value = '0.425'
$ws.Range("E25").Value = '  +5.25%  '

$ws.Range("D26").Value = '2.795.26'
$ws.Range("E26").Value = '  +4.70%  '

$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '0.162'
$ws.Range("E27").Value = '  +3.19%  '

$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '0.990'
$ws.Range("E28").Value = '  -0.94%  '

$ws.Range("D29").Value = '0.0₃0825'
$ws.Range("E29").Value = '  +7.01%  '

$ws.Range("D30").Value = '7.28'
$ws.Range("E30").Value = '  +6.60%  '

$ws.Range("D31").Value = '0.996'
$ws.Range("E31").Value = '  -0.47%  '

$ws.Range("D32").Value = '6.46'
$ws.Range("E32").Value = '  +12.11%  '

$ws.Range("D33").Value = '19.20'
$ws.Range("E33").Value = '  +4.03%  '

$ws.Range("D34").Value = '1.59'
$ws.Range("E34").Value = '  +3.69%  '

$ws.Range("D35").Value = '150.13'
$ws.Range("E35").Value = '  +0.62%  '

$ws.Range("D36").Value = '1.03'
$ws.Range("E36").Value = '  +18.22%  '

$ws.Range("D37").Value = '4.10'
$ws.Range("E37").Value = '  +6.05%  '

$ws.Range("D38").Value = '1.16'
$ws.Range("E38").Value = '  +5.63%  '

$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").Value = '0.861'
$ws.Range("E39").Value = '  +5.26%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '37.00'
$ws.Range("E40").Value = '  +3.51%  '

$ws.Range("D41").Value = '3.73'
$ws.Range("E41").Value = '  +7.28%  '

$ws.Range("E42").Value = '  +3.42%  '

$ws.Range("D43").Value = '0.629'
$ws.Range("E43").Value = '  +4.28%  '

$ws.Range("D44").Value = '283.59'
$ws.Range("E44").Value = '  -0.42%  '

$ws.Range("D45").Value = '20.36'
$ws.Range("E45").Value = '  +9.95%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '0.0986'
$ws.Range("E46").Value = '  +0.84%  '

$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").Value = '0.993'
$ws.Range("E47").Value = '  -0.48%  '

$ws.Range("D48").Value = '0.0538'
$ws.Range("E48").Value = '  +2.07%  '

$ws.Range("D49").Value = '4.80'
$ws.Range("E49").Value = '  +6.50%  '

$ws.Range("D50").Value = '2.026.74'
$ws.Range("E50").Value = '  +7.54%  '

$ws.Range("D51").Value = '0.0232'
$ws.Range("E51").Value = '  +3.31%  '

